$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the updated "through" date
$ws.Name = "Through 2021-10-18"

# Update the month label for October to reflect the new "through" date
$ws.Range("A11").Value = "October (through 10-18)"

# Update H8 (July 2021) value
$ws.Range("H8").Value = 150

# Update row 11 (October) values for years 2015-2021
$ws.Range("B11").Value = 17
$ws.Range("C11").Value = 29
$ws.Range("D11").Value = 31
$ws.Range("E11").Value = 46
$ws.Range("F11").Value = 28
$ws.Range("G11").Value = 86
$ws.Range("H11").Value = 112

# Update row 12 (Total) values for years 2015-2021
$ws.Range("B12").Value = 243
$ws.Range("C12").Value = 458
$ws.Range("D12").Value = 658
$ws.Range("E12").Value = 594
$ws.Range("F12").Value = 450
$ws.Range("G12").Value = 987
$ws.Range("H12").Value = 1360
